# Correct 2229v1 growth rates on the RateCompare sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RateCompare")

# Tspan for strain 2229v1 (row 2) was mis-entered as "2-8"; correct value is "10-30".
$ws.Range("B2").Value = "10-30"

# Corrected growth-rate figures for the same row.
$ws.Range("C2").Value = 0.13
$ws.Range("D2").Value = 0.03
$ws.Range("E2").Value = 1.53
$ws.Range("F2").Value = 0.86

# Move the active selection to F3, matching the saved view state.
[void]$ws.Range("F3").Select()
